$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2321.4285
$ws.Range("I40").Value = 2742.8572
$ws.Range("J40").Value = 1900
$ws.Range("K40").Value = 2742.8572
$ws.Range("L40").Value = 1900
$ws.Range("M40").Value = -2567.8572
$ws.Range("N40").Value = -2250

$ws.Range("H100").Value = 2559.2727
$ws.Range("I100").Value = 2706.5
$ws.Range("J100").Value = 2166.6667
$ws.Range("K100").Value = 2706.5
$ws.Range("L100").Value = 2166.6667
$ws.Range("M100").Value = -2165.5
$ws.Range("N100").Value = -3248.6667

$ws.Range("H106").Value = 10184.857
$ws.Range("I106").Value = 10506.308
$ws.Range("K106").Value = 10506.308
$ws.Range("M106").Value = -9875.308000000001

$ws.Range("H111").Value = 4166.125
$ws.Range("I111").Value = 3765.8
$ws.Range("J111").Value = 4833.3335
$ws.Range("K111").Value = 11297.4
$ws.Range("L111").Value = 14500.0005
$ws.Range("M111").Value = -8230.400000000001
$ws.Range("N111").Value = -20634.0005

$ws.Range("H116").Value = 2027.35
$ws.Range("I116").Value = 1636.0667
$ws.Range("K116").Value = 1636.0667
$ws.Range("M116").Value = 1805.9333

$ws.Range("H129").Value = 650.9048
$ws.Range("I129").Value = 344.0909
$ws.Range("J129").Value = 988.4
$ws.Range("K129").Value = 1032.2727
$ws.Range("L129").Value = 2965.2
$ws.Range("M129").Value = 3967.7273
$ws.Range("N129").Value = -12965.2

$ws.Range("H132").Value = 8553738
$ws.Range("I132").Value = 9264329
$ws.Range("J132").Value = 26635.334
$ws.Range("K132").Value = 27792987
$ws.Range("L132").Value = 79906.00199999999
$ws.Range("M132").Value = -27790457
$ws.Range("N132").Value = -84966.00199999999

$ws.Range("H141").Value = 949.8889
$ws.Range("I141").Value = 693.625
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 2080.875
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 3099.125
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4195.041
$ws.Range("I32").Value = 3809.7812
$ws.Range("J32").Value = 6934.6665
$ws.Range("K32").Value = 3809.7812
$ws.Range("L32").Value = 6934.6665
$ws.Range("M32").Value = -3522.7812
$ws.Range("N32").Value = -7508.6665

$ws.Range("H45").Value = 1533.1111
$ws.Range("I45").Value = 1571.1428
$ws.Range("K45").Value = 1571.1428
$ws.Range("M45").Value = -1194.1428

$ws.Range("H74").Value = 952.8958
$ws.Range("I74").Value = 738.375
$ws.Range("J74").Value = 2025.5
$ws.Range("K74").Value = 738.375
$ws.Range("L74").Value = 2025.5
$ws.Range("M74").Value = 135.625
$ws.Range("N74").Value = -3773.5

$ws.Range("H77").Value = 952.8958
$ws.Range("I77").Value = 738.375
$ws.Range("J77").Value = 2025.5
$ws.Range("K77").Value = 3691.875
$ws.Range("L77").Value = 10127.5
$ws.Range("M77").Value = 676.125
$ws.Range("N77").Value = -18863.5

$ws.Range("H110").Value = 2411
$ws.Range("I110").Value = 2100
$ws.Range("J110").Value = 2473.2
$ws.Range("K110").Value = 2100
$ws.Range("L110").Value = 2473.2
$ws.Range("M110").Value = -55
$ws.Range("N110").Value = -6563.2

$ws.Range("H122").Value = 2387.5454
$ws.Range("I122").Value = 2658.2222
$ws.Range("J122").Value = 1169.5
$ws.Range("K122").Value = 7974.6666
$ws.Range("L122").Value = 3508.5
$ws.Range("M122").Value = -5524.6666
$ws.Range("N122").Value = -8408.5

$ws.Range("H132").Value = 2356.2144
$ws.Range("I132").Value = 2513.0476
$ws.Range("K132").Value = 7539.1428
$ws.Range("M132").Value = -5009.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2033.909
$ws.Range("I20").Value = 1609.2307
$ws.Range("J20").Value = 2647.3333
$ws.Range("K20").Value = 1609.2307
$ws.Range("L20").Value = 2647.3333
$ws.Range("M20").Value = -1362.2307
$ws.Range("N20").Value = -3141.3333

$ws.Range("H134").Value = 4068.8647
$ws.Range("I134").Value = 1078.0646
$ws.Range("K134").Value = 3234.1938
$ws.Range("M134").Value = -699.1938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2504.55
$ws.Range("I31").Value = 2692.7334
$ws.Range("K31").Value = 2692.7334
$ws.Range("M31").Value = -2397.7334

$ws.Range("H34").Value = 2504.55
$ws.Range("I34").Value = 2692.7334
$ws.Range("K34").Value = 2692.7334
$ws.Range("M34").Value = -2490.7334

$ws.Range("H62").Value = 6062982
$ws.Range("I62").Value = 2380
$ws.Range("J62").Value = 66669000
$ws.Range("K62").Value = 2380
$ws.Range("L62").Value = 66669000
$ws.Range("M62").Value = -1756
$ws.Range("N62").Value = -66670248

$ws.Range("H65").Value = 6062982
$ws.Range("I65").Value = 2380
$ws.Range("J65").Value = 66669000
$ws.Range("K65").Value = 11900
$ws.Range("L65").Value = 333345000
$ws.Range("M65").Value = -8780
$ws.Range("N65").Value = -333351240

$ws.Range("H121").Value = 8463
$ws.Range("J121").Value = 8463
$ws.Range("L121").Value = 8463
$ws.Range("N121").Value = -11083

$ws.Range("H132").Value = 4224.05
$ws.Range("I132").Value = 4468.6763
$ws.Range("J132").Value = 2837.8333
$ws.Range("K132").Value = 13406.0289
$ws.Range("L132").Value = 8513.499899999999
$ws.Range("M132").Value = -10876.0289
$ws.Range("N132").Value = -13573.4999

$ws.Range("H141").Value = 29328.076
$ws.Range("J141").Value = 29328.076
$ws.Range("L141").Value = 29328.076
$ws.Range("N141").Value = -39688.076

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9909.944
$ws.Range("I3").Value = 3998.5715
$ws.Range("J3").Value = 13671.728
$ws.Range("K3").Value = 11995.7145
$ws.Range("L3").Value = 41015.18399999999
$ws.Range("M3").Value = -11883.7145
$ws.Range("N3").Value = -41239.18399999999

$ws.Range("H13").Value = 499.33334
$ws.Range("I13").Value = 249.5
$ws.Range("K13").Value = 748.5
$ws.Range("M13").Value = -580.5

$ws.Range("H81").Value = 2616.9
$ws.Range("J81").Value = 3027.5334
$ws.Range("L81").Value = 9082.600199999999
$ws.Range("N81").Value = -11328.6002

$ws.Range("H84").Value = 2616.9
$ws.Range("J84").Value = 3027.5334
$ws.Range("L84").Value = 27247.8006
$ws.Range("N84").Value = -38479.8006

$ws.Range("H105").Value = 107099.9
$ws.Range("J105").Value = 107099.9
$ws.Range("L105").Value = 321299.7
$ws.Range("N105").Value = -326541.7

$ws.Range("H106").Value = 2559.7273
$ws.Range("J106").Value = 2715.7
$ws.Range("L106").Value = 8147.099999999999
$ws.Range("N106").Value = -10039.1

$ws.Range("H107").Value = 5631.4736
$ws.Range("J107").Value = 6998.6
$ws.Range("L107").Value = 20995.8
$ws.Range("N107").Value = -24835.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2026.5
$ws.Range("I102").Value = 1377.3158
$ws.Range("J102").Value = 6138
$ws.Range("K102").Value = 1377.3158
$ws.Range("L102").Value = 6138
$ws.Range("M102").Value = 244.6841999999999
$ws.Range("N102").Value = -9382

$ws.Range("H113").Value = 1535.1538
$ws.Range("I113").Value = 1449.3334
$ws.Range("K113").Value = 1449.3334
$ws.Range("M113").Value = 720.6666

$ws.Range("H132").Value = 2707.8
$ws.Range("I132").Value = 2342.2666
$ws.Range("J132").Value = 3804.4
$ws.Range("K132").Value = 7026.7998
$ws.Range("L132").Value = 11413.2
$ws.Range("M132").Value = -4496.7998
$ws.Range("N132").Value = -16473.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 832.0526
$ws.Range("I22").Value = 520.6667
$ws.Range("K22").Value = 520.6667
$ws.Range("M22").Value = -225.6667

$ws.Range("H27").Value = 832.0526
$ws.Range("I27").Value = 520.6667
$ws.Range("K27").Value = 520.6667
$ws.Range("M27").Value = -413.6667

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()

$ws.Range("H40").Value = 4008.2727
$ws.Range("I40").Value = 2454
$ws.Range("J40").Value = 11002.5
$ws.Range("K40").Value = 2454
$ws.Range("L40").Value = 11002.5
$ws.Range("M40").Value = -2318
$ws.Range("N40").Value = -11274.5

$ws.Range("H46").Value = 4090.923
$ws.Range("I46").Value = 1095
$ws.Range("J46").Value = 5422.4443
$ws.Range("K46").Value = 1095
$ws.Range("L46").Value = 5422.4443
$ws.Range("M46").Value = -907
$ws.Range("N46").Value = -5798.4443

$ws.Range("H136").Value = 2561.2034
$ws.Range("I136").Value = 2559.125
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 7677.375
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -5127.375
$ws.Range("N136").Value = -12900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 264.96295
$ws.Range("I113").Value = 164.90909
$ws.Range("J113").Value = 333.75
$ws.Range("K113").Value = 494.72727
$ws.Range("L113").Value = 1001.25
$ws.Range("M113").Value = 1675.27273
$ws.Range("N113").Value = -5341.25

$ws.Range("H126").Value = 38463036
$ws.Range("I126").Value = 50001236
$ws.Range("J126").Value = 2366.5
$ws.Range("K126").Value = 150003708
$ws.Range("L126").Value = 7099.5
$ws.Range("M126").Value = -150001238
$ws.Range("N126").Value = -12039.5

$ws.Range("H132").Value = 5606.2
$ws.Range("I132").Value = 8393.666999999999
$ws.Range("K132").Value = 25181.001
$ws.Range("M132").Value = -22651.001
